$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure Price column cells that receive text values (even numeric-looking ones) stay as text
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D51").NumberFormat = "@"

$ws.Range("D2").Value = '98.519.52'
$ws.Range("E2").Value = '  +4.81%  '

$ws.Range("D3").Value = '3.366.06'
$ws.Range("E3").Value = '  +9.66%  '

$ws.Range("E4").Value = '  +0.00%  '

$ws.Range("D5").Value = '256.61'
$ws.Range("E5").Value = '  +8.73%  '

$ws.Range("D6").Value = '623.23'
$ws.Range("E6").Value = '  +2.59%  '

$ws.Range("E7").Value = '  +8.22%  '

$ws.Range("D8").Value = '0.386'
$ws.Range("E8").Value = '  +1.81%  '

$ws.Range("E9").Value = '  +0.02%  '

$ws.Range("D10").Value = '3.361.50'
$ws.Range("E10").Value = '  +9.64%  '

$ws.Range("D11").Value = '0.812'
$ws.Range("E11").Value = '  +0.96%  '

$ws.Range("E12").Value = '  +1.70%  '

$ws.Range("D13").Value = '98.188.57'
$ws.Range("E13").Value = '  +4.82%  '

$ws.Range("D14").Value = '35.70'
$ws.Range("E14").Value = '  +5.87%  '

$ws.Range("D15").Value = '0.0000247'
$ws.Range("E15").Value = '  +2.73%  '

$ws.Range("D16").Value = '3.997.06'
$ws.Range("E16").Value = '  +9.83%  '

$ws.Range("D17").Value = '5.49'
$ws.Range("E17").Value = '  +3.56%  '

$ws.Range("D18").Value = '3.367.57'
$ws.Range("E18").Value = '  +9.47%  '

$ws.Range("D19").Value = '3.67'
$ws.Range("E19").Value = '  +3.30%  '

$ws.Range("D20").Value = '15.03'
$ws.Range("E20").Value = '  +5.02%  '

$ws.Range("D21").Value = '485.82'
$ws.Range("E21").Value = '  +10.02%  '

$ws.Range("E22").Value = '  +3.15%  '

$ws.Range("D23").Value = '0.0000208'
$ws.Range("E23").Value = '  +10.12%  '

$ws.Range("D24").Value = '9.25'
$ws.Range("E24").Value = '  +4.93%  '

$ws.Range("D25").Value = '5.72'
$ws.Range("E25").Value = '  +3.83%  '

$ws.Range("D26").Value = '87.98'
$ws.Range("E26").Value = '  +4.12%  '

$ws.Range("D27").Value = '12.05'
$ws.Range("E27").Value = '  +1.40%  '

$ws.Range("E28").Value = '  +9.69%  '

$ws.Range("E29").Value = '  +0.03%  '

$ws.Range("D30").Value = '0.252'
$ws.Range("E30").Value = '  +0.96%  '

$ws.Range("D31").Value = '0.186'
$ws.Range("E31").Value = '  +4.26%  '

$ws.Range("D32").Value = '0.127'
$ws.Range("E32").Value = '  +2.95%  '

$ws.Range("D33").Value = '1.00'
$ws.Range("E33").Value = '  +0.12%  '

$ws.Range("D34").Value = '9.24'
$ws.Range("E34").Value = '  +3.81%  '

$ws.Range("D35").Value = '27.45'
$ws.Range("E35").Value = '  +8.54%  '

$ws.Range("D36").Value = '523.53'
$ws.Range("E36").Value = '  +7.53%  '

$ws.Range("E37").Value = '  +0.05%  '

$ws.Range("D38").Value = '7.36'
$ws.Range("E38").Value = '  -0.83%  '

$ws.Range("E39").Value = '  +4.08%  '

$ws.Range("D40").Value = '24.81'
$ws.Range("E40").Value = '  +3.26%  '

$ws.Range("D41").Value = '0.449'
$ws.Range("E41").Value = '  +3.27%  '

$ws.Range("E42").Value = '  -2.38%  '

$ws.Range("E43").Value = '  +2.98%  '

$ws.Range("D44").Value = '3.25'
$ws.Range("E44").Value = '  +5.68%  '

$ws.Range("D45").Value = '0.783'
$ws.Range("E45").Value = '  +16.48%  '

$ws.Range("E46").Value = '  -0.02%  '

$ws.Range("D47").Value = '160.96'
$ws.Range("E47").Value = '  -0.22%  '

$ws.Range("E48").Value = '  +6.57%  '

$ws.Range("B49").Value = 'Filecoin'
$ws.Range("C49").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("D49").Value = '4.54'
$ws.Range("E49").Value = '  +6.84%  '

$ws.Range("B50").Value = 'OKB'
$ws.Range("C50").Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range("D50").Value = '45.41'
$ws.Range("E50").Value = '  +4.31%  '

$ws.Range("B51").Value = 'ImmutableX'
$ws.Range("C51").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D51").Value = '1.36'
$ws.Range("E51").Value = '  +6.43%  '
